$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 31; this shifts existing rows 31-49 down to 32-50
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new weekly data point
$ws.Cells.Item(31, 1).Value = 3
$ws.Cells.Item(31, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(31, 3).Value = "Coquimbo"
$ws.Cells.Item(31, 4).Value = 44567
$ws.Cells.Item(31, 5).Value = 5
$ws.Cells.Item(31, 6).Value = 100112022
$ws.Cells.Item(31, 7).Value = "Arveja Verde"
$ws.Cells.Item(31, 8).Value = "Perfection"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 68
$ws.Cells.Item(31, 11).Value = 24000
$ws.Cells.Item(31, 12).Value = 25000
$ws.Cells.Item(31, 13).Value = 24559
$ws.Cells.Item(31, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(31, 16).Value = 982
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
